$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all target cells to Text format to preserve exact string representation
# (prevents Excel from auto-converting numeric-looking strings like '0.597' or
# dates/number-looking strings like '63.937.30' into numbers).
$targetCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "E9", "D10", "E10", "D11", "E11", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "E35", "D36", "E36", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "E41", "E42", "B43", "C43", "D43", "E43", "B44", "C44", "D44", "E44", "D45", "E45", "B46", "C46", "D46", "E46", "B47", "C47", "D47", "E47", "B48", "C48", "D48", "E48", "D49", "D50", "E50", "D51", "E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cryptocurrency data values
$ws.Range('D2').Value = '63.937.30'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '3.344.12'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '535.59'
$ws.Range('E5').Value = '  +3.15%  '
$ws.Range('D6').Value = '174.74'
$ws.Range('E6').Value = '  -5.39%  '
$ws.Range('D7').Value = '0.597'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.341.44'
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').Value = '0.612'
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('D11').Value = '53.95'
$ws.Range('E11').Value = '  -10.80%  '
$ws.Range('E12').Value = '  +2.11%  '
$ws.Range('D13').Value = '0.0000259'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').Value = '9.26'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').Value = '3.883.97'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '3.349.00'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('D18').Value = '17.56'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = '63.596.01'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  +1.73%  '
$ws.Range('D21').Value = '0.969'
$ws.Range('E21').Value = '  +1.65%  '
$ws.Range('D22').Value = '373.04'
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('D23').Value = '4.18'
$ws.Range('E23').Value = '  +6.95%  '
$ws.Range('D24').Value = '11.37'
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('D25').Value = '3.78'
$ws.Range('E25').Value = '  +2.33%  '
$ws.Range('D26').Value = '81.68'
$ws.Range('E26').Value = '  +1.27%  '
$ws.Range('D27').Value = '6.19'
$ws.Range('E27').Value = '  +4.47%  '
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('D29').Value = '11.35'
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('D30').Value = '8.32'
$ws.Range('E30').Value = '  -0.92%  '
$ws.Range('D31').Value = '29.01'
$ws.Range('E31').Value = '  +1.24%  '
$ws.Range('D32').Value = '648.96'
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('D33').Value = '6.52'
$ws.Range('E33').Value = '  -3.91%  '
$ws.Range('D34').Value = '11.28'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  +1.37%  '
$ws.Range('D36').Value = '58.63'
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').Value = '37.20'
$ws.Range('E38').Value = '  +1.86%  '
$ws.Range('D39').Value = '0.385'
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('D40').Value = '0.0₃0731'
$ws.Range('E40').Value = '  +10.92%  '
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.934.92'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '2.59'
$ws.Range('E44').Value = '  +6.58%  '
$ws.Range('D45').Value = '2.98'
$ws.Range('E45').Value = '  +4.51%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0401'
$ws.Range('E46').Value = '  +2.49%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '3.14'
$ws.Range('E47').Value = '  +5.36%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '2.66'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('D49').Value = '2.63'
$ws.Range('D50').Value = '0.126'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').Value = '137.76'
$ws.Range('E51').Value = '  +4.75%  '
